$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.15625130904473
$ws.Range("D2").Value = 5.227209632459942
$ws.Range("E2").Value = 18.47924742512029
$ws.Range("F2").Value = 25.93875600156668
$ws.Range("G2").Value = 3.641485923773893
$ws.Range("K2").Value = 10.15479932006108
$ws.Range("L2").Value = 8.318050084480998
$ws.Range("O2").Value = 23.20648075587856
$ws.Range("B3").Value = 16.04732137356886
$ws.Range("D3").Value = 5.185903575125639
$ws.Range("E3").Value = 18.55333116928583
$ws.Range("F3").Value = 25.94978019416371
$ws.Range("G3").Value = 3.64343243343379
$ws.Range("K3").Value = 9.842001205313155
$ws.Range("L3").Value = 8.272406162376789
$ws.Range("O3").Value = 23.26565112086701
$ws.Range("B4").Value = 15.98340499647301
$ws.Range("D4").Value = 5.160001923818734
$ws.Range("E4").Value = 18.60168733642252
$ws.Range("F4").Value = 25.96371436656088
$ws.Range("G4").Value = 3.644691410686243
$ws.Range("K4").Value = 9.643019092030556
$ws.Range("L4").Value = 8.245675200544103
$ws.Range("O4").Value = 23.30735015492503
$ws.Range("B5").Value = 15.95812820486309
$ws.Range("D5").Value = 5.149315180079005
$ws.Range("E5").Value = 18.62211498678112
$ws.Range("F5").Value = 25.97119250243004
$ws.Range("G5").Value = 3.645220550912595
$ws.Range("K5").Value = 9.56026665272573
$ws.Range("L5").Value = 8.235116209079349
$ws.Range("O5").Value = 23.32568930681319
$ws.Range("B6").Value = 15.95397816187143
$ws.Range("D6").Value = 5.147532825084801
$ws.Range("E6").Value = 18.62555062139207
$ws.Range("F6").Value = 25.97254286966011
$ws.Range("G6").Value = 3.645309387938821
$ws.Range("K6").Value = 9.546427411653207
$ws.Range("L6").Value = 8.233383315987604
$ws.Range("O6").Value = 23.32881572681916
$ws.Range("B7").Value = 15.98306095833745
$ws.Range("D7").Value = 5.159858326150003
$ws.Range("E7").Value = 18.60195990591186
$ws.Range("F7").Value = 25.96380793518942
$ws.Range("G7").Value = 3.644698481613583
$ws.Range("K7").Value = 9.64190970312084
$ws.Range("L7").Value = 8.245531434598806
$ws.Range("O7").Value = 23.30759203611693
$ws.Range("B8").Value = 16.11809096645482
$ws.Range("D8").Value = 5.213081714077939
$ws.Range("E8").Value = 18.5041966491363
$ws.Range("F8").Value = 25.9410691651007
$ws.Range("G8").Value = 3.642143866546568
$ws.Range("K8").Value = 10.04842626551504
$ws.Range("L8").Value = 8.302048636883283
$ws.Range("O8").Value = 23.22576654840323
$ws.Range("B9").Value = 16.40533534314651
$ws.Range("D9").Value = 5.312996639051848
$ws.Range("E9").Value = 18.33520552743734
$ws.Range("F9").Value = 25.95336996646668
$ws.Range("G9").Value = 3.63763828782224
$ws.Range("K9").Value = 10.78772965064303
$ws.Range("L9").Value = 8.422754870550813
$ws.Range("O9").Value = 23.10804429214461
$ws.Range("B10").Value = 16.6285276361887
$ws.Range("D10").Value = 5.383459509149121
$ws.Range("E10").Value = 18.2248475831291
$ws.Range("F10").Value = 25.99707228999403
$ws.Range("G10").Value = 3.634632066691105
$ws.Range("K10").Value = 11.2921466241652
$ws.Range("L10").Value = 8.516902096493212
$ws.Range("O10").Value = 23.04778840345215
$ws.Range("B11").Value = 16.73235870765641
$ws.Range("D11").Value = 5.414828367107423
$ws.Range("E11").Value = 18.17762866547568
$ws.Range("F11").Value = 26.02445204789939
$ws.Range("G11").Value = 3.633329788045239
$ws.Range("K11").Value = 11.51259757925401
$ws.Range("L11").Value = 8.560792990475118
$ws.Range("O11").Value = 23.02610580495455
$ws.Range("B12").Value = 16.77197767542121
$ws.Range("D12").Value = 5.426604504094419
$ws.Range("E12").Value = 18.1601763832777
$ws.Range("F12").Value = 26.03589391183795
$ws.Range("G12").Value = 3.632845982139759
$ws.Range("K12").Value = 11.5947376582873
$ws.Range("L12").Value = 8.577555142323471
$ws.Range("O12").Value = 23.01872100448346
$ws.Range("B13").Value = 16.76343214864302
$ws.Range("D13").Value = 5.424072927771715
$ws.Range("E13").Value = 18.16391599502251
$ws.Range("F13").Value = 26.03338203414908
$ws.Range("G13").Value = 3.632949763790772
$ws.Range("K13").Value = 11.5771075635484
$ws.Range("L13").Value = 8.573938995113279
$ws.Range("O13").Value = 23.02027468807912
$ws.Range("B14").Value = 16.73561233869989
$ws.Range("D14").Value = 5.415799275558616
$ws.Range("E14").Value = 18.17618427089189
$ws.Range("F14").Value = 26.02537188863395
$ws.Range("G14").Value = 3.633289798144781
$ws.Range("K14").Value = 11.51938238129792
$ws.Range("L14").Value = 8.56216924945605
$ws.Range("O14").Value = 23.02548168855856
$ws.Range("B15").Value = 16.71861011535545
$ws.Range("D15").Value = 5.410717948205782
$ws.Range("E15").Value = 18.18375473203263
$ws.Range("F15").Value = 26.02060511676674
$ws.Range("G15").Value = 3.633499293813136
$ws.Range("K15").Value = 11.48384832260759
$ws.Range("L15").Value = 8.554978047366161
$ws.Range("O15").Value = 23.0287787474607
$ws.Range("B16").Value = 16.62178572787007
$ws.Range("D16").Value = 5.381395346579788
$ws.Range("E16").Value = 18.22799344023202
$ws.Range("F16").Value = 25.99543349548475
$ws.Range("G16").Value = 3.634718483008874
$ws.Range("K16").Value = 11.27755441868328
$ws.Range("L16").Value = 8.514054157146736
$ws.Range("O16").Value = 23.04932087523303
$ws.Range("B17").Value = 16.56295470983273
$ws.Range("D17").Value = 5.363228506996314
$ws.Range("E17").Value = 18.25589621267972
$ws.Range("F17").Value = 25.98190937381387
$ws.Range("G17").Value = 3.635483099200191
$ws.Range("K17").Value = 11.14865878305282
$ws.Range("L17").Value = 8.489212991767955
$ws.Range("O17").Value = 23.06339160907891
$ws.Range("B18").Value = 16.52933497726494
$ws.Range("D18").Value = 5.352715178344789
$ws.Range("E18").Value = 18.27222599680876
$ws.Range("F18").Value = 25.97483688971976
$ws.Range("G18").Value = 3.635929032149596
$ws.Range("K18").Value = 11.07367556794317
$ws.Range("L18").Value = 8.475025737706892
$ws.Range("O18").Value = 23.07202376730125
$ws.Range("B19").Value = 16.5179903251545
$ws.Range("D19").Value = 5.349144637153873
$ws.Range("E19").Value = 18.27780323324999
$ws.Range("F19").Value = 25.97256368552538
$ws.Range("G19").Value = 3.636081074377713
$ws.Range("K19").Value = 11.04814369579103
$ws.Range("L19").Value = 8.470239815194208
$ws.Range("O19").Value = 23.07503897548288
$ws.Range("B20").Value = 16.56919499574627
$ws.Range("D20").Value = 5.365169076654539
$ws.Range("E20").Value = 18.25289685212343
$ws.Range("F20").Value = 25.9832759768803
$ws.Range("G20").Value = 3.635401068757151
$ws.Range("K20").Value = 11.16246782285627
$ws.Range("L20").Value = 8.4918470318626
$ws.Range("O20").Value = 23.061837949085
$ws.Range("B21").Value = 16.74377578221727
$ws.Range("D21").Value = 5.418232263375359
$ws.Range("E21").Value = 18.17256915562367
$ws.Range("F21").Value = 26.02769556476503
$ws.Range("G21").Value = 3.633189668729798
$ws.Range("K21").Value = 11.53637435759756
$ws.Range("L21").Value = 8.565622555886614
$ws.Range("O21").Value = 23.02392983472748
$ws.Range("B22").Value = 16.85961103072511
$ws.Range("D22").Value = 5.452312453420963
$ws.Range("E22").Value = 18.1225677475933
$ws.Range("F22").Value = 26.06298188895318
$ws.Range("G22").Value = 3.631798803077711
$ws.Range("K22").Value = 11.77291695718627
$ws.Range("L22").Value = 8.614659269181166
$ws.Range("O22").Value = 23.00396962744428
$ws.Range("B23").Value = 16.79763879060515
$ws.Range("D23").Value = 5.434179402861633
$ws.Range("E23").Value = 18.14902609489137
$ws.Range("F23").Value = 26.04357839201204
$ws.Range("G23").Value = 3.632536170801711
$ws.Range("K23").Value = 11.64739914804217
$ws.Range("L23").Value = 8.588416166997261
$ws.Range("O23").Value = 23.0141815502957
$ws.Range("B24").Value = 16.5663731278452
$ws.Range("D24").Value = 5.364291959624675
$ws.Range("E24").Value = 18.25425196555951
$ws.Range("F24").Value = 25.98265594637
$ws.Range("G24").Value = 3.635438134953037
$ws.Range("K24").Value = 11.15622749172597
$ws.Range("L24").Value = 8.490655887873677
$ws.Range("O24").Value = 23.06253866841863
$ws.Range("B25").Value = 16.32537919972787
$ws.Range("D25").Value = 5.286467826032693
$ws.Range("E25").Value = 18.37849502074108
$ws.Range("F25").Value = 25.94394868693312
$ws.Range("G25").Value = 3.638803543653114
$ws.Range("K25").Value = 10.59430730127924
$ws.Range("L25").Value = 8.389099304551898
$ws.Range("O25").Value = 23.13529611861646
